$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: change from "borrar" node to "Importación" node; drop D value (posicionY)
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "Importación"
$ws.Range("C44").Value = 0
$ws.Range("D44").ClearContents()

# Row 45: change from "borrar" node to "Exportación" node; drop D value (posicionY)
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "Exportación"
$ws.Range("C45").Value = 4
$ws.Range("D45").ClearContents()

# Row 46: keep as "borrar" but renumber id
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = "borrar"
$ws.Range("C46").Value = 6
$ws.Range("D46").Value = 1

# Row 47: keep as "borrar" but renumber id
$ws.Range("A47").Value = 45
$ws.Range("B47").Value = "borrar"
$ws.Range("C47").Value = 6
$ws.Range("D47").Value = 1

# Row 48: keep as "borrar" but renumber id
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = "borrar"
$ws.Range("C48").Value = 6
$ws.Range("D48").Value = 1

# Update view to reflect scrolled/selected state from the edit session
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("C51").Select()
